$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.874784666666667
$ws.Range("N2").Value = 8.624354
$ws.Range("O2").Value = 0.1187109652550681
$ws.Range("P2").Value = 0.121184727686443
$ws.Range("Q2").Value = 1.225395171136222
$ws.Range("R2").Value = 11.028556540226
$ws.Range("S2").Value = 0.1187109652550681
$ws.Range("T2").Value = 0.121184727686443

# Row 3
$ws.Range("O3").Value = 0.4442422727481699
$ws.Range("P3").Value = 0.4534996302499962
$ws.Range("S3").Value = 0.4442422727481699
$ws.Range("T3").Value = 0.4534996302499962

# Row 4
$ws.Range("M4").Value = 5.147441999999999
$ws.Range("N4").Value = 15.442326
$ws.Range("O4").Value = 0.2125577666737049
$ws.Range("P4").Value = 0.2169871588243338
$ws.Range("Q4").Value = 2.194129752966
$ws.Range("R4").Value = 19.747167776694
$ws.Range("S4").Value = 0.2125577666737049
$ws.Range("T4").Value = 0.2169871588243338

# Row 5
$ws.Range("M5").Value = 1.483016
$ws.Range("N5").Value = 2.966032
$ws.Range("O5").Value = 0.06123946008548931
$ws.Range("P5").Value = 0.04167706708575228
$ws.Range("Q5").Value = 0.6321449624346668
$ws.Range("R5").Value = 3.792869774608
$ws.Range("S5").Value = 0.06123946008548931
$ws.Range("T5").Value = 0.04167706708575228

# Row 6
$ws.Range("M6").Value = 3.953360666666667
$ws.Range("N6").Value = 11.860082
$ws.Range("O6").Value = 0.1632495352375677
$ws.Range("P6").Value = 0.1666514161534747
$ws.Range("Q6").Value = 1.685145022117556
$ws.Range("R6").Value = 15.166305199058
$ws.Range("S6").Value = 0.1632495352375677
$ws.Range("T6").Value = 0.1666514161534747
